# Add three new "coding question" entries (rows 72-74) to the tracker sheet,
# matching the author's "Add files via upload" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 72: "24. Swap Nodes in Pairs in LinkedList." -----------------------
$ws.Range("B72").Value = 45767
$ws.Range("C72").Value = "24. Swap Nodes in Pairs in LinkedList. Leetcode"
$ws.Range("C72").Characters(40, 8).Font.Bold = $true
$ws.Range("F72").Value = "O(n)"
$ws.Range("G72").Value = "solved and submitted"

# --- Row 73: "Merge k Sorted Lists." -----------------------------------------
$ws.Range("B73").Value = 45767
$ws.Range("C73").Value = "Merge k Sorted Lists. Leetcode"
$ws.Range("C73").Characters(23, 8).Font.Bold = $true
$ws.Range("F73").Value = "solved using merge sort. Better approach can be piority queue"

# --- Row 74: "Rotate List by k." ---------------------------------------------
$ws.Range("B74").Value = 45767
$ws.Range("C74").Value = "Rotate List by k. Leetcode"
$ws.Range("C74").Characters(19, 8).Font.Bold = $true
$ws.Range("G74").Value = "solved and submitted"

# Date column keeps the sheet's existing dd/mm/yyyy custom display format,
# picked up automatically from column B's style.

# Reflect where the author was last working when they saved.
[void]$ws.Range("B74").Select()
